$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.245.27"
$ws.Range("E2").Value = "  +0.46%  "

$ws.Range("D3").Value = "3.597.73"
$ws.Range("E3").Value = "  +2.45%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.996"
$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "605.54"
$ws.Range("E5").Value = "  -0.40%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.46"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.625"
$ws.Range("E7").Value = "  -0.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.206"
$ws.Range("E9").Value = "  -2.83%  "

$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.97"
$ws.Range("E11").Value = "  -0.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000305"
$ws.Range("E12").Value = "  -0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.56"
$ws.Range("E13").Value = "  -0.83%  "

$ws.Range("D14").Value = "4.151.19"
$ws.Range("E14").Value = "  +2.10%  "

$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.09"
$ws.Range("E15").Value = "  +2.86%  "

$ws.Range("B16").Value = "BitcoinCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "600.42"
$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.29"
$ws.Range("E17").Value = "  +1.20%  "

$ws.Range("D18").Value = "70.264.01"
$ws.Range("E18").Value = "  +0.41%  "

$ws.Range("D19").Value = "3.569.21"
$ws.Range("E19").Value = "  +1.92%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.123"
$ws.Range("E20").Value = "  +1.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.996"
$ws.Range("E21").Value = "  -0.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.89"
$ws.Range("E22").Value = "  +0.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "103.04"
$ws.Range("E23").Value = "  -0.86%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.15"
$ws.Range("E24").Value = "  +1.04%  "

$ws.Range("E25").Value = "  -0.51%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.10"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.86"
$ws.Range("E27").Value = "  -1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.61"
$ws.Range("E28").Value = "  -2.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.06"
$ws.Range("E29").Value = "  +0.66%  "

$ws.Range("E30").Value = "  -1.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.30"
$ws.Range("E31").Value = "  -6.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.39"
$ws.Range("E32").Value = "  -3.43%  "

$ws.Range("E33").Value = "  +0.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.58"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("D35").Value = "3.904.23"
$ws.Range("E35").Value = "  +5.09%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.25"
$ws.Range("E36").Value = "  +7.95%  "

$ws.Range("D37").Value = "0.0₃0828"
$ws.Range("E37").Value = "  +2.89%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "528.28"
$ws.Range("E38").Value = "  +0.57%  "

$ws.Range("E39").Value = "  +0.15%  "

$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.20"
$ws.Range("E40").Value = "  +0.69%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  +1.30%  "

$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.394"
$ws.Range("E42").Value = "  +0.21%  "

$ws.Range("E43").Value = "  -2.46%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0455"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.87"
$ws.Range("E45").Value = "  -0.29%  "

$ws.Range("E46").Value = "  +0.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.32"
$ws.Range("E47").Value = "  +0.31%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").Value = "  -2.14%  "

$ws.Range("E49").Value = "  -0.01%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000251"
$ws.Range("E50").Value = "  +4.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.31"
$ws.Range("E51").Value = "  +2.98%  "
